# Add new row 6 to Sheet1: variable "e" with a "normal" distribution from numpy.random,
# param1 = 1, and a computed param (E6/10) - supports the new 'single_var' sensitivity
# analysis flag described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = "e"
$ws.Range("C6").Value = "numpy.random"
$ws.Range("D6").Value = "normal"
$ws.Range("E6").Value = 1
$ws.Range("F6").Formula = "=E6/10"

# Match the formatting already used by the other rows in the table (C/D columns).
$ws.Range("C6").Font.Color = 0
$ws.Range("D6").Font.Color = 0

# Move the active selection down to E11 (below the newly inserted data row).
$ws.Range("E11").Select()
